$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Preserve the empty "quote-prefix" style (F column, row that will become row 2) ---
# F7 currently carries style s=8 (vertical-top, quotePrefix) and no value. After the
# reorder, that same blank/styled cell lands in F2. Copy its formatting there first,
# then clear any value row 2 might already have in that column.
$ws.Range("F7").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("F2").ClearContents()

# --- Rewrite the data rows (A2:F9) in the new order / with the new labels ---
# Row 2: "Anabelle's" method, now first, relabeled A / A1 / A2 (the old "4-year"
# sub-row is dropped).
$ws.Range("A2").Value = "A. Lowest consecutive flows:in Reclamation's ensembles:and traces (2025)"
$ws.Range("B2").Value = "Lee Ferry"
$ws.Range("C2").Value = "Natural"
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0

$ws.Range("A3").Value = "A1. 10-year"
$ws.Range("B3").Value = "Lee Ferry"
$ws.Range("C3").Value = "Natural"
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = 7.5
$ws.Range("F3").Value = "Immersive modeling in progress"

$ws.Range("A4").Value = "A2. 3-year"
$ws.Range("B4").Value = "Lee Ferry"
$ws.Range("C4").Value = "Natural"
$ws.Range("D4").Value = 3.9
$ws.Range("E4").Value = 6.8
$ws.Range("F4").Value = "Immersive modeling in progress"

$ws.Range("A5").Value = "B. From tree rings back to:1400 AD (2023)"
$ws.Range("B5").Value = "Lee Ferry"
$ws.Range("C5").Value = "Natural"
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 14
$ws.Range("F5").Value = "Cap depletions to 10-year lookback:period of flow."

$ws.Range("A6").Value = "C. Collaborator choices in:immersive modeling:sessions (2021)"
$ws.Range("B6").Value = "Lee Ferry"
$ws.Range("C6").Value = "Natural"
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = 13
$ws.Range("F6").Value = "Divide inflow; Users consume and:conserve within their:account balance."

$ws.Range("A7").Value = "D. 85%, 65%, and 50% of:2000 to 2018 average:flow (2022)"
$ws.Range("B7").Value = "Lake Powell"
$ws.Range("C7").Value = "Regulated Inflow"
$ws.Range("D7").Value = 6
$ws.Range("E7").Value = 10
$ws.Range("F7").Value = "Release 95% of regulated:inflow."

$ws.Range("A8").Value = "E. Release to prevent:drawdown to 3,490 feet:(2024)."
$ws.Range("B8").Value = "Lake Powell"
$ws.Range("C8").Value = "Release"
$ws.Range("D8").Value = 4
$ws.Range("E8").Value = 6
$ws.Range("F8").Value = "Release to prevent drawdown:to 3,490 feet."

$ws.Range("A9").Value = "F. Low Lake Powell:releases + gains through:Grand Canyon (2022)"
$ws.Range("B9").Value = "Lake Mead"
$ws.Range("C9").Value = "Regulated Inflow"
$ws.Range("D9").Value = 7
$ws.Range("E9").Value = 10
$ws.Range("F9").Value = "Rule curve; Consumption equals or:less than inflow minus:evaporation."

# The old table had 9 data rows (2-10); the new one has only 8 (2-9), so drop the
# now-unused trailing row.
$ws.Rows.Item(10).Delete()

# --- View state: zoom + selection as left by the edit ---
$ws.Application.ActiveWindow.Zoom = 150
$ws.Range("A10").Select()
